$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new expense entry is being inserted at row 15 ("[SMAC] Venus", 550),
# pushing the "Tong" (total) row that used to live at row 15 down to row 18.

# 1) Move the total row's formatting from row 15 down to row 18.
$ws.Range("G15:I15").Copy()
$ws.Range("G18:I18").PasteSpecial(-4122)

# 2) Re-create the total row's content at its new location (row 18),
#    now summing through the new row 15 entry.
$ws.Range("G18").Value = "Tổng"
$ws.Range("H18").Formula = "=SUM(H3:H15)"
$ws.Range("I18").Formula = "=SUM(I3:I14)"

# 3) Give row 15 the plain data-row formatting (same as other entry rows,
#    e.g. row 8) and drop the now-unused I15 cell entirely.
$ws.Range("G8").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("G8").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("I15").Clear()

# 4) Fill in the new expense entry itself.
$ws.Range("G15").Value = "[SMAC] Venus"
$ws.Range("H15").Value = 550

# 5) Match the author's final selection.
[void]$ws.Range("G16").Select()
